# Updated cryptos list (Price / Volume(1h) columns) to match the latest scrape.
# Values are written as text (matching the source workbook's inlineStr cells) by
# forcing a Text number format before the write, then resetting the style back to
# "Normal" so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.871.77"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +4.08%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.878.71"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.07%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "278.33"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.03%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5309"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +4.12%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3442"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.26%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "45.15"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.06962"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.40%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.09"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.8047"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -3.05%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.07742"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.54%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.883.83"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +4.35%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "90.38"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.37%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "5.175"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.95%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.56"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.07%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008034"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "26.928.29"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +4.08%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.121.41"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("E23").Value = "  +0.68%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "10.04"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "6.210"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  +7.70%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "147.11"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +4.28%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.665"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "17.35"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "113.68"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.91%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.353"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.316"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.08891"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.06%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.04908"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +3.41%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7271"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.890"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.291"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.60%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.358"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.10%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01850"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.5117"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.56%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9564"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "115.86"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.81%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "6.189"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.106"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  -0.01%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4469"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.1343"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.309"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "36.22"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05952"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.94%  "
